$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a cell so it is stored as literal TEXT (shared string),
# even when its content looks like a number (e.g. "212756 " with a
# trailing space). Using a formula that evaluates to a string and then
# pasting-by-value avoids Excel's automatic text->number coercion while
# keeping the destination cell's existing number format / style (so no
# new style entries get created in styles.xml).
# ---------------------------------------------------------------------
function Set-TextValue($range, [string]$text) {
    $helper = $ws.Range("ZZ1")
    $escaped = $text.Replace("""", """""")
    $helper.Formula = "=""" + $escaped + """"
    $helper.Copy()
    $range.PasteSpecial(-4163)
    $helper.Clear()
}

# Row 110: BOMBILLAS LED SMART | 2755 (number)
$ws.Range("A109").Copy()
$ws.Range("A110").PasteSpecial(-4122)
$ws.Range("A110").Value = "BOMBILLAS LED SMART"

$ws.Range("C109").Copy()
$ws.Range("C110").PasteSpecial(-4122)
$ws.Range("C110").Value = 2755

# Row 111: BOMBILLAS LED SMART | "212756 " (text)
$ws.Range("A109").Copy()
$ws.Range("A111").PasteSpecial(-4122)
$ws.Range("A111").Value = "BOMBILLAS LED SMART"

$ws.Range("C109").Copy()
$ws.Range("C111").PasteSpecial(-4122)
Set-TextValue $ws.Range("C111") "212756 "

# Row 112: BOMBILLAS LED SMART | "3001 " (text)
$ws.Range("A109").Copy()
$ws.Range("A112").PasteSpecial(-4122)
$ws.Range("A112").Value = "BOMBILLAS LED SMART"

$ws.Range("C109").Copy()
$ws.Range("C112").PasteSpecial(-4122)
Set-TextValue $ws.Range("C112") "3001 "

# Row 113: BOMBILLAS LED SMART | "2999 " (text)
$ws.Range("A109").Copy()
$ws.Range("A113").PasteSpecial(-4122)
$ws.Range("A113").Value = "BOMBILLAS LED SMART"

$ws.Range("C109").Copy()
$ws.Range("C113").PasteSpecial(-4122)
Set-TextValue $ws.Range("C113") "2999 "

# Row 114: BOMBILLAS LED SMART | "212754 " (text)
$ws.Range("A109").Copy()
$ws.Range("A114").PasteSpecial(-4122)
$ws.Range("A114").Value = "BOMBILLAS LED SMART"

$ws.Range("C109").Copy()
$ws.Range("C114").PasteSpecial(-4122)
Set-TextValue $ws.Range("C114") "212754 "

# Row 115: BOMBILLAS LED SMART | "2755 " (text)
$ws.Range("A109").Copy()
$ws.Range("A115").PasteSpecial(-4122)
$ws.Range("A115").Value = "BOMBILLAS LED SMART"

$ws.Range("C109").Copy()
$ws.Range("C115").PasteSpecial(-4122)
Set-TextValue $ws.Range("C115") "2755 "

# Row 116: BOMBILLAS LED SMART | "2751 " (text)
$ws.Range("A109").Copy()
$ws.Range("A116").PasteSpecial(-4122)
$ws.Range("A116").Value = "BOMBILLAS LED SMART"

$ws.Range("C109").Copy()
$ws.Range("C116").PasteSpecial(-4122)
Set-TextValue $ws.Range("C116") "2751 "

# Row 117: APARATOS SMART | 8445 (number)
$ws.Range("A109").Copy()
$ws.Range("A117").PasteSpecial(-4122)
$ws.Range("A117").Value = "APARATOS SMART"

$ws.Range("C109").Copy()
$ws.Range("C117").PasteSpecial(-4122)
$ws.Range("C117").Value = 8445

$excel.CutCopyMode = 0

# Update the view to match the final selection/scroll position.
$ws.Range("B100").Select()

Write-Output "done"
